$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.337.96"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").Value = "1.871.49"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.26"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4715"
$ws.Range("E7").Value = "  +0.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2895"
$ws.Range("E8").Value = "  +2.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06642"
$ws.Range("E9").Value = "  +1.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.65"
$ws.Range("E10").Value = "  -0.30%  "

$ws.Range("E11").Value = "  +1.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.52"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").Value = "1.871.32"
$ws.Range("E13").Value = "  +0.12%  "

$ws.Range("E14").Value = "  -0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6887"
$ws.Range("E15").Value = "  +1.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "272.18"
$ws.Range("E16").Value = "  -2.68%  "

$ws.Range("D17").Value = "30.324.09"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.22"
$ws.Range("E18").Value = "  +6.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007780"
$ws.Range("E19").Value = "  +6.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.09%  "

$ws.Range("D21").Value = "2.117.10"
$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.331"
$ws.Range("E22").Value = "  -1.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.222"
$ws.Range("E24").Value = "  +0.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.21"
$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.355"
$ws.Range("E26").Value = "  +1.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.00"
$ws.Range("E27").Value = "  -0.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.959"
$ws.Range("E28").Value = "  +1.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.375"
$ws.Range("E29").Value = "  -0.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09981"
$ws.Range("E30").Value = "  +2.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.381"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("E33").Value = "  +0.49%  "

$ws.Range("E34").Value = "  -0.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  +0.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7029"
$ws.Range("E36").Value = "  -0.48%  "

$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01886"
$ws.Range("E38").Value = "  +1.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.654"
$ws.Range("E39").Value = "  +2.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.323"
$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.83"
$ws.Range("E41").Value = "  -2.85%  "

$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8443"
$ws.Range("E43").Value = "  -0.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4173"
$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.27"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.258"
$ws.Range("E47").Value = "  -0.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.115"
$ws.Range("E48").Value = "  -1.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "937.91"
$ws.Range("E49").Value = "  -2.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.58"
$ws.Range("E50").Value = "  +1.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05677"
$ws.Range("E51").Value = "  +0.53%  "
